$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.103.16'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '1.779.19'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '329.13'
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '0.4509'
$ws.Range("E7").Value = '  +1.48%  '
$ws.Range("D8").Value = '0.3566'
$ws.Range("E8").Value = '  +1.21%  '
$ws.Range("D9").Value = '0.07456'
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").Value = '42.07'
$ws.Range("E10").Value = '  +1.32%  '
$ws.Range("D11").Value = '1.109'
$ws.Range("E11").Value = '  +2.91%  '
$ws.Range("D12").Value = '0.9997'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '21.05'
$ws.Range("E13").Value = '  +2.91%  '
$ws.Range("D14").Value = '6.063'
$ws.Range("E14").Value = '  +2.73%  '
$ws.Range("D15").Value = '7.276'
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("D16").Value = '1.774.44'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").Value = '94.05'
$ws.Range("E17").Value = '  +2.87%  '
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("D19").Value = '0.06451'
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").Value = '0.9991'
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("D22").Value = '5.805'
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").Value = '28.086.00'
$ws.Range("E23").Value = '  +1.93%  '
$ws.Range("E24").Value = '  +1.94%  '
$ws.Range("D25").Value = '2.127'
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").Value = '161.77'
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("E27").Value = '  +1.88%  '
$ws.Range("D28").Value = '1.979.34'
$ws.Range("E28").Value = '  +1.98%  '
$ws.Range("D29").Value = '2.169'
$ws.Range("E29").Value = '  +6.63%  '
$ws.Range("D30").Value = '125.11'
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").Value = '1.110'
$ws.Range("E31").Value = '  +6.03%  '
$ws.Range("D32").Value = '5.732'
$ws.Range("E32").Value = '  +6.77%  '
$ws.Range("D33").Value = '0.09219'
$ws.Range("E33").Value = '  +1.51%  '
$ws.Range("D34").Value = '3.697'
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").Value = '11.91'
$ws.Range("E35").Value = '  +2.71%  '
$ws.Range("D36").Value = '0.06215'
$ws.Range("E36").Value = '  +3.19%  '
$ws.Range("D37").Value = '0.02297'
$ws.Range("E37").Value = '  +1.12%  '
$ws.Range("D38").Value = '0.2116'
$ws.Range("E38").Value = '  +2.73%  '
$ws.Range("D39").Value = '5.006'
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("D40").Value = '0.6334'
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").Value = '1.397'
$ws.Range("E42").Value = '  +1.88%  '
$ws.Range("D43").Value = '7.937'
$ws.Range("E43").Value = '  +2.95%  '
$ws.Range("D44").Value = '13.28'
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("E45").Value = '  +1.48%  '
$ws.Range("E46").Value = '  +2.34%  '
$ws.Range("D47").Value = '122.90'
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("E48").Value = '  +2.19%  '
$ws.Range("E49").Value = '  +2.71%  '
$ws.Range("D50").Value = '0.06902'
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("D51").Value = '73.04'
$ws.Range("E51").Value = '  +2.38%  '